$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily format the data range as text so that numeric-looking strings
# (e.g. "1.00", "2.99") are not auto-converted to numbers when assigned.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "43.953.11"
$ws.Range("E2").Value = "  -4.76%  "
$ws.Range("D3").Value = "2.660.88"
$ws.Range("E3").Value = "  +2.78%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "305.47"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").Value = "97.13"
$ws.Range("E6").Value = "  -1.52%  "
$ws.Range("D7").Value = "0.589"
$ws.Range("E7").Value = "  -1.41%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").Value = "0.567"
$ws.Range("E9").Value = "  -1.49%  "
$ws.Range("D10").Value = "37.68"
$ws.Range("E10").Value = "  -2.99%  "
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("E12").Value = "  -2.53%  "
$ws.Range("D13").Value = "3.069.44"
$ws.Range("E13").Value = "  +3.05%  "
$ws.Range("E14").Value = "  +1.84%  "
$ws.Range("D15").Value = "2.677.55"
$ws.Range("E15").Value = "  +3.70%  "
$ws.Range("D16").Value = "0.906"
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("D17").Value = "14.81"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").Value = "44.100.70"
$ws.Range("E18").Value = "  -4.70%  "
$ws.Range("D19").Value = "6.83"
$ws.Range("E19").Value = "  +3.04%  "
$ws.Range("D20").Value = "0.0₃0991"
$ws.Range("E20").Value = "  -1.25%  "
$ws.Range("D21").Value = "12.61"
$ws.Range("E21").Value = "  -2.47%  "
$ws.Range("D22").Value = "74.18"
$ws.Range("E22").Value = "  +3.70%  "
$ws.Range("D23").Value = "273.61"
$ws.Range("E23").Value = "  +1.12%  "
$ws.Range("E24").Value = "  +5.94%  "
$ws.Range("D25").Value = "2.99"
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").Value = "30.28"
$ws.Range("E26").Value = "  +1.13%  "
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "10.41"
$ws.Range("E28").Value = "  -0.64%  "
$ws.Range("D29").Value = "2.24"
$ws.Range("E29").Value = "  -2.83%  "
$ws.Range("D30").Value = "37.64"
$ws.Range("E30").Value = "  -3.49%  "
$ws.Range("E31").Value = "  -0.78%  "
$ws.Range("D32").Value = "3.66"
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("E33").Value = "  +6.71%  "
$ws.Range("D34").Value = "153.51"
$ws.Range("E34").Value = "  +2.80%  "
$ws.Range("E35").Value = "  -1.77%  "
$ws.Range("D36").Value = "0.0827"
$ws.Range("E36").Value = "  -1.19%  "
$ws.Range("E37").Value = "  -1.39%  "
$ws.Range("D38").Value = "25.09"
$ws.Range("E38").Value = "  +5.90%  "
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("D40").Value = "15.87"
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("E42").Value = "  -1.89%  "
$ws.Range("D43").Value = "3.91"
$ws.Range("E43").Value = "  -4.06%  "
$ws.Range("D44").Value = "2.117.59"
$ws.Range("E44").Value = "  -0.58%  "
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("D46").Value = "90.98"
$ws.Range("E46").Value = "  -2.91%  "
$ws.Range("D47").Value = "9.22"
$ws.Range("E47").Value = "  -4.04%  "
$ws.Range("D48").Value = "2.924.47"
$ws.Range("E48").Value = "  +3.23%  "
$ws.Range("D49").Value = "108.88"
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("E50").Value = "  +3.83%  "
$ws.Range("E51").Value = "  -1.61%  "

# Restore the original (default) cell formatting now that the text values are set.
$dataRange.ClearFormats()

